$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns starting at column K (11) to make room for the new
# "Flash SW" fields. This shifts all existing columns from K onward (K..GD)
# to the right by 4 (becoming O..GH).
$ws.Range("K1:N1").EntireColumn.Insert()

# Fill in the headers for the newly inserted columns in row 1
$ws.Range("K1").Value = "BCMe Flash SW"
$ws.Range("L1").Value = "CMC1 Flash SW"
$ws.Range("M1").Value = "CMC2 Flash SW"
$ws.Range("N1").Value = "CMC3 Flash SW"

# Fill in the corresponding data values for row 2 (stored as text, matching
# the other cells in this sheet which are all inline strings)
$ws.Range("K2:N2").NumberFormat = "@"
$ws.Range("K2").Value = "1000"
$ws.Range("L2").Value = "1000"
$ws.Range("M2").Value = "1000"
$ws.Range("N2").Value = "1000"
